# Update Name of Algo
# Apply numeric corrections to result_data_RandomForest.xlsx (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = 5.977999999999993
$ws.Range("B4").Value  = 5.105100000000006
$ws.Range("D6").Value  = -7.963799999999993
$ws.Range("B7").Value  = 5.733199999999993
$ws.Range("D7").Value  = -7.26209999999999
$ws.Range("B8").Value  = 5.617399999999996
$ws.Range("D8").Value  = -7.497799999999996
$ws.Range("A11").Value = -21.79700000000001
$ws.Range("A12").Value = -22.69400000000001
$ws.Range("B12").Value = 5.154500000000001
$ws.Range("B14").Value = 8.871100000000002
$ws.Range("A15").Value = -21.44860000000002
$ws.Range("D19").Value = -8.58649999999999
$ws.Range("D21").Value = -7.702200000000002
$ws.Range("B22").Value = 5.538300000000001
$ws.Range("D24").Value = -7.566599999999998
$ws.Range("D25").Value = -7.910599999999991
